$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Append a new row to Sheet1, re-using the formatting of the row above it,
# carrying over the "kumar" record that used to sit in Sheet2!A6.
$ws1.Range("A4").Copy($ws1.Range("A5"))
$ws1.Range("A5").Value = "kumar"
$ws1.Range("B5").Value = 3

# Clear the value that used to live in Sheet2!A6 (keep its formatting/style)
$ws2.Range("A6").ClearContents()

# Update the recorded selection on each sheet
$ws1.Range("A5").Select()
$ws2.Activate()
$ws2.Range("A5").Select()
